$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 514.3333
$ws.Range("I12").Value = 477.6
$ws.Range("J12").Value = 698
$ws.Range("K12").Value = 477.6
$ws.Range("L12").Value = 698
$ws.Range("M12").Value = -307.6
$ws.Range("N12").Value = -1038
$ws.Range("H18").Value = 2588.4167
$ws.Range("I18").Value = 2732.6365
$ws.Range("K18").Value = 2732.6365
$ws.Range("M18").Value = -2448.6365
$ws.Range("H32").Value = 1942.1364
$ws.Range("I32").Value = 1716.4286
$ws.Range("J32").Value = 2047.4667
$ws.Range("K32").Value = 1716.4286
$ws.Range("L32").Value = 2047.4667
$ws.Range("M32").Value = -1390.4286
$ws.Range("N32").Value = -2699.4667
$ws.Range("H43").Value = 30305360
$ws.Range("J43").Value = 2702.25
$ws.Range("L43").Value = 2702.25
$ws.Range("N43").Value = -2840.25
$ws.Range("H52").Value = 2084.7273
$ws.Range("I52").Value = 2284.5715
$ws.Range("J52").Value = 1735
$ws.Range("K52").Value = 6853.7145
$ws.Range("L52").Value = 5205
$ws.Range("M52").Value = -6693.7145
$ws.Range("N52").Value = -5525
$ws.Range("H57").Value = 27316.5
$ws.Range("J57").Value = 27316.5
$ws.Range("L57").Value = 81949.5
$ws.Range("N57").Value = -82947.5
$ws.Range("H58").Value = 1785
$ws.Range("J58").Value = 2333.3333
$ws.Range("L58").Value = 6999.999899999999
$ws.Range("N58").Value = -7299.999899999999
$ws.Range("H98").Value = 1856.3334
$ws.Range("I98").Value = 1931.5883
$ws.Range("J98").Value = 1536.5
$ws.Range("K98").Value = 1931.5883
$ws.Range("L98").Value = 1536.5
$ws.Range("M98").Value = -433.5882999999999
$ws.Range("N98").Value = -4532.5
$ws.Range("H111").Value = 1722.3572
$ws.Range("I111").Value = 1338.6
$ws.Range("K111").Value = 4015.8
$ws.Range("M111").Value = -948.7999999999997
$ws.Range("H122").Value = 1856.3334
$ws.Range("I122").Value = 1931.5883
$ws.Range("J122").Value = 1536.5
$ws.Range("K122").Value = 5794.7649
$ws.Range("L122").Value = 4609.5
$ws.Range("M122").Value = -3344.7649
$ws.Range("N122").Value = -9509.5
$ws.Range("H127").Value = 3000
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H129").Value = 2002.8235
$ws.Range("I129").Value = 1548.5
$ws.Range("J129").Value = 2063.4
$ws.Range("K129").Value = 4645.5
$ws.Range("L129").Value = 6190.200000000001
$ws.Range("M129").Value = 354.5
$ws.Range("N129").Value = -16190.2
$ws.Range("H132").Value = 19423.95
$ws.Range("I132").Value = 24999.205
$ws.Range("J132").Value = 1901.7142
$ws.Range("K132").Value = 74997.61500000001
$ws.Range("L132").Value = 5705.142599999999
$ws.Range("M132").Value = -72467.61500000001
$ws.Range("N132").Value = -10765.1426
$ws.Range("H137").Value = 1798.8889
$ws.Range("I137").Value = 1282
$ws.Range("K137").Value = 3846
$ws.Range("M137").Value = -1296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 319484.12
$ws.Range("I32").Value = 373388.25
$ws.Range("K32").Value = 373388.25
$ws.Range("M32").Value = -373101.25
$ws.Range("H102").Value = 2205.389
$ws.Range("I102").Value = 2376.8125
$ws.Range("K102").Value = 2376.8125
$ws.Range("M102").Value = -754.8125
$ws.Range("H122").Value = 2491.9473
$ws.Range("I122").Value = 2197.4119
$ws.Range("K122").Value = 6592.2357
$ws.Range("M122").Value = -4142.2357
$ws.Range("H132").Value = 3771.3333
$ws.Range("I132").Value = 1822
$ws.Range("K132").Value = 5466
$ws.Range("M132").Value = -2936

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 1190
$ws.Range("I23").Value = 1190
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1190
$ws.Range("N23").ClearContents()
$ws.Range("M23").Value = -907
$ws.Range("H105").Value = 13345.833
$ws.Range("I105").Value = 15395.25
$ws.Range("K105").Value = 15395.25
$ws.Range("M105").Value = -13648.25
$ws.Range("H107").Value = 8225.378000000001
$ws.Range("I107").Value = 10974.52
$ws.Range("K107").Value = 10974.52
$ws.Range("M107").Value = -9054.52
$ws.Range("H134").Value = 3652.6562
$ws.Range("I134").Value = 2275.2307
$ws.Range("K134").Value = 6825.6921
$ws.Range("M134").Value = -4290.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3040.2222
$ws.Range("I16").Value = 1681.6875
$ws.Range("K16").Value = 1681.6875
$ws.Range("M16").Value = -1394.6875
$ws.Range("H31").Value = 2397.8704
$ws.Range("I31").Value = 2130.5417
$ws.Range("J31").Value = 2611.7334
$ws.Range("K31").Value = 2130.5417
$ws.Range("L31").Value = 2611.7334
$ws.Range("M31").Value = -1835.5417
$ws.Range("N31").Value = -3201.7334
$ws.Range("H34").Value = 2397.8704
$ws.Range("I34").Value = 2130.5417
$ws.Range("J34").Value = 2611.7334
$ws.Range("K34").Value = 2130.5417
$ws.Range("L34").Value = 2611.7334
$ws.Range("M34").Value = -1928.5417
$ws.Range("N34").Value = -3015.7334
$ws.Range("H58").Value = 1902.5714
$ws.Range("I58").Value = 1453.9166
$ws.Range("K58").Value = 1453.9166
$ws.Range("M58").Value = -1250.9166
$ws.Range("H86").Value = 11504.556
$ws.Range("J86").Value = 12379.9
$ws.Range("L86").Value = 12379.9
$ws.Range("N86").Value = -14625.9
$ws.Range("H89").Value = 11504.556
$ws.Range("J89").Value = 12379.9
$ws.Range("L89").Value = 61899.5
$ws.Range("N89").Value = -73131.5
$ws.Range("H99").Value = 3822.2222
$ws.Range("J99").Value = 4023.5557
$ws.Range("L99").Value = 4023.5557
$ws.Range("N99").Value = -7019.5557
$ws.Range("H103").Value = 19812.2
$ws.Range("I103").Value = 19812.2
$ws.Range("K103").Value = 19812.2
$ws.Range("M103").Value = -18640.2
$ws.Range("H107").Value = 1727.0869
$ws.Range("I107").Value = 1515.6
$ws.Range("J107").Value = 2123.625
$ws.Range("K107").Value = 1515.6
$ws.Range("L107").Value = 2123.625
$ws.Range("M107").Value = 404.4000000000001
$ws.Range("N107").Value = -5963.625
$ws.Range("H113").Value = 3040.2222
$ws.Range("I113").Value = 1681.6875
$ws.Range("K113").Value = 1681.6875
$ws.Range("M113").Value = 488.3125
$ws.Range("H123").Value = 169990
$ws.Range("J123").Value = 169990
$ws.Range("L123").Value = 169990
$ws.Range("N123").Value = -179790
$ws.Range("H126").Value = 3822.2222
$ws.Range("J126").Value = 4023.5557
$ws.Range("L126").Value = 12070.6671
$ws.Range("N126").Value = -17010.6671
$ws.Range("H134").Value = 2301.276
$ws.Range("I134").Value = 1671
$ws.Range("J134").Value = 2633
$ws.Range("K134").Value = 5013
$ws.Range("L134").Value = 7899
$ws.Range("M134").Value = -2478
$ws.Range("N134").Value = -12969
$ws.Range("H136").Value = 1902.5714
$ws.Range("I136").Value = 1453.9166
$ws.Range("K136").Value = 4361.7498
$ws.Range("M136").Value = -1811.7498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 378.55554
$ws.Range("I2").Value = 201.14285
$ws.Range("J2").Value = 999.5
$ws.Range("K2").Value = 1206.8571
$ws.Range("L2").Value = 5997
$ws.Range("M2").Value = -1093.8571
$ws.Range("N2").Value = -6223
$ws.Range("H12").Value = 972.2941
$ws.Range("J12").Value = 955.4
$ws.Range("L12").Value = 2866.2
$ws.Range("N12").Value = -3212.2
$ws.Range("H38").Value = 785.4211
$ws.Range("J38").Value = 1257.3
$ws.Range("L38").Value = 3771.9
$ws.Range("N38").Value = -4465.9
$ws.Range("H39").Value = 4951.6113
$ws.Range("J39").Value = 5320.5625
$ws.Range("L39").Value = 15961.6875
$ws.Range("N39").Value = -16549.6875
$ws.Range("H55").Value = 4969.533
$ws.Range("J55").Value = 5832.9165
$ws.Range("L55").Value = 17498.7495
$ws.Range("N55").Value = -17852.7495
$ws.Range("H56").Value = 7999.125
$ws.Range("I56").Value = 7999.125
$ws.Range("K56").Value = 7999.125
$ws.Range("M56").Value = -7469.125
$ws.Range("H62").Value = 14125.571
$ws.Range("I62").Value = 6250
$ws.Range("J62").Value = 17275.8
$ws.Range("K62").Value = 18750
$ws.Range("L62").Value = 51827.39999999999
$ws.Range("M62").Value = -18064
$ws.Range("N62").Value = -53199.39999999999
$ws.Range("H65").Value = 14125.571
$ws.Range("I65").Value = 6250
$ws.Range("J65").Value = 17275.8
$ws.Range("K65").Value = 56250
$ws.Range("L65").Value = 155482.2
$ws.Range("M65").Value = -52818
$ws.Range("N65").Value = -162346.2
$ws.Range("H70").Value = 8541.857
$ws.Range("I70").Value = 4499.5
$ws.Range("K70").Value = 13498.5
$ws.Range("M70").Value = -13183.5
$ws.Range("H73").Value = 8541.857
$ws.Range("I73").Value = 4499.5
$ws.Range("K73").Value = 13498.5
$ws.Range("M73").Value = -12406.5
$ws.Range("H98").Value = 627.2857
$ws.Range("J98").Value = 678.6
$ws.Range("L98").Value = 2035.8
$ws.Range("N98").Value = -5031.8
$ws.Range("H131").Value = 8360128.5
$ws.Range("I131").Value = 12990077
$ws.Range("J131").Value = 257718.75
$ws.Range("K131").Value = 38970231
$ws.Range("L131").Value = 773156.25
$ws.Range("M131").Value = -38965191
$ws.Range("N131").Value = -783236.25
$ws.Range("H132").Value = 1563.7
$ws.Range("I132").Value = 1555.2858
$ws.Range("J132").Value = 1583.3334
$ws.Range("K132").Value = 13997.5722
$ws.Range("L132").Value = 14250.0006
$ws.Range("M132").Value = -11467.5722
$ws.Range("N132").Value = -19310.0006
$ws.Range("H138").Value = 4021.9583
$ws.Range("I138").Value = 3346.625
$ws.Range("K138").Value = 10039.875
$ws.Range("M138").Value = -4899.875
$ws.Range("H139").Value = 5684916.5
$ws.Range("I139").Value = 9617075
$ws.Range("K139").Value = 28851225
$ws.Range("M139").Value = -28846085

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 121.73333
$ws.Range("I2").Value = 54
$ws.Range("J2").Value = 166.88889
$ws.Range("K2").Value = 54
$ws.Range("L2").Value = 166.88889
$ws.Range("M2").Value = 59
$ws.Range("N2").Value = -392.88889
$ws.Range("H12").Value = 500
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H80").Value = 93961.266
$ws.Range("I80").Value = 129750.55
$ws.Range("K80").Value = 129750.55
$ws.Range("M80").Value = -128752.55
$ws.Range("H83").Value = 93961.266
$ws.Range("I83").Value = 129750.55
$ws.Range("K83").Value = 648752.75
$ws.Range("M83").Value = -643760.75
$ws.Range("H113").Value = 2391.913
$ws.Range("I113").Value = 2221.4707
$ws.Range("K113").Value = 2221.4707
$ws.Range("M113").Value = -51.47069999999985

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 4899.5
$ws.Range("I30").Value = 3866
$ws.Range("K30").Value = 3866
$ws.Range("M30").Value = -3758
$ws.Range("H40").Value = 1144902.9
$ws.Range("I40").Value = 1514914.8
$ws.Range("K40").Value = 1514914.8
$ws.Range("M40").Value = -1514778.8
$ws.Range("H100").Value = 4282.4
$ws.Range("J100").Value = 4378.5
$ws.Range("L100").Value = 4378.5
$ws.Range("N100").Value = -5460.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 39999
$ws.Range("J68").Value = 39999
$ws.Range("L68").Value = 39999
$ws.Range("N68").Value = -41621
$ws.Range("H71").Value = 39999
$ws.Range("J71").Value = 39999
$ws.Range("L71").Value = 119997
$ws.Range("N71").Value = -128109
$ws.Range("H113").Value = 498
$ws.Range("J113").Value = 697.3333
$ws.Range("L113").Value = 2091.9999
$ws.Range("N113").Value = -6431.9999
$ws.Range("H122").Value = 14162.917
$ws.Range("I122").Value = 16036.477
$ws.Range("J122").Value = 1048
$ws.Range("K122").Value = 48109.431
$ws.Range("L122").Value = 3144
$ws.Range("M122").Value = -45659.431
$ws.Range("N122").Value = -8044
$ws.Range("H124").Value = 29872
$ws.Range("J124").Value = 29872
$ws.Range("L124").Value = 29872
$ws.Range("N124").Value = -39692
$ws.Range("H132").Value = 25857.342
$ws.Range("I132").Value = 38239.926
$ws.Range("K132").Value = 114719.778
$ws.Range("M132").Value = -112189.778
$ws.Range("H136").Value = 33404.395
$ws.Range("I136").Value = 60186.53
$ws.Range("J136").Value = 4948.375
$ws.Range("K136").Value = 180559.59
$ws.Range("L136").Value = 14845.125
$ws.Range("M136").Value = -178009.59
$ws.Range("N136").Value = -19945.125
